$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update row 11 values
$ws.Range("B11").Value = 3
$ws.Range("C11").Value = 1.5
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.5
$ws.Range("I11").Value = "Team Project setup"

# Update selection
$ws.Range("E13").Select()
